# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# OFF sheet - Row 3 (Road "R" totals)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 190
$wsOff.Range("C3").Value = 131
$wsOff.Range("D3").Value = 41
$wsOff.Range("E3").Value = 25

# DEF sheet - Row 3 (Road "R" totals)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 225
$wsDef.Range("C3").Value = 168
$wsDef.Range("D3").Value = 64
$wsDef.Range("E3").Value = 27
$wsDef.Range("F3").Value = 2
$wsDef.Range("G3").Value = 4
